$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values: the source text often has trailing zeros or
# multi-dot thousand separators (e.g. "0.9980", "26.783.67") that must be
# preserved verbatim as text, not coerced into a floating point number.
# Force the cell to Text format before writing, then clear the format again
# afterwards so the cell keeps its original (default/general) style.
$dValues = [ordered]@{
    "D2" = '26.783.67';
    "D3" = '1.727.09';
    "D4" = '0.9973';
    "D5" = '240.63';
    "D6" = '0.9980';
    "D7" = '0.4809';
    "D9" = '0.06179';
    "D10" = '1.721.80';
    "D11" = '15.85';
    "D13" = '0.6034';
    "D14" = '4.463';
    "D15" = '76.90';
    "D16" = '0.9981';
    "D17" = '26.581.39';
    "D18" = '0.9976';
    "D19" = '0.000007123';
    "D20" = '11.37';
    "D21" = '1.945.37';
    "D22" = '4.411';
    "D23" = '8.479';
    "D24" = '5.062';
    "D25" = '139.91';
    "D27" = '1.781';
    "D28" = '106.46';
    "D30" = '3.979';
    "D31" = '0.07922';
    "D32" = '3.670';
    "D33" = '0.04526';
    "D35" = '0.9998';
    "D36" = '0.6184';
    "D37" = '0.9294';
    "D38" = '1.996';
    "D39" = '2.443';
    "D40" = '0.9972';
    "D41" = '0.01494';
    "D42" = '5.613';
    "D43" = '99.72';
    "D44" = '0.3829';
    "D45" = '6.776';
    "D47" = '0.05357';
    "D48" = '7.899';
    "D51" = '51.51';
}

foreach ($addr in $dValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$addr]
    $cell.ClearFormats()
}

# Columns B, C and E are plain text already (coin names, links, and
# whitespace-padded "+/-x.xx%" strings), so they can be written directly.
$otherValues = [ordered]@{
    "E3" = '  +0.24%  ';
    "E4" = '  -0.09%  ';
    "E5" = '  -0.56%  ';
    "E6" = '  -0.08%  ';
    "E7" = '  -1.58%  ';
    "E8" = '  -0.35%  ';
    "E9" = '  -0.07%  ';
    "E10" = '  -0.06%  ';
    "E11" = '  +2.16%  ';
    "E12" = '  -1.96%  ';
    "E13" = '  +0.88%  ';
    "E14" = '  -1.11%  ';
    "E15" = '  -0.16%  ';
    "E16" = '  -0.08%  ';
    "E17" = '  +0.71%  ';
    "E18" = '  -0.08%  ';
    "E19" = '  -0.21%  ';
    "E20" = '  +0.74%  ';
    "E21" = '  -0.09%  ';
    "E22" = '  -0.72%  ';
    "E23" = '  -0.22%  ';
    "E25" = '  +1.33%  ';
    "E26" = '  -0.07%  ';
    "E27" = '  +2.77%  ';
    "E28" = '  +0.18%  ';
    "E29" = '  -2.05%  ';
    "E30" = '  +1.94%  ';
    "E31" = '  -1.37%  ';
    "E32" = '  +0.61%  ';
    "E33" = '  +0.78%  ';
    "E34" = '  -0.34%  ';
    "E35" = '  +0.40%  ';
    "E36" = '  -0.45%  ';
    "E37" = '  +0.98%  ';
    "B38" = 'RenderToken';
    "C38" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr';
    "E38" = '  +1.82%  ';
    "B39" = 'MXToken';
    "C39" = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';
    "E39" = '  +2.63%  ';
    "E40" = '  -0.09%  ';
    "E41" = '  +1.16%  ';
    "E42" = '  +3.19%  ';
    "E43" = '  -0.17%  ';
    "E44" = '  -0.26%  ';
    "E45" = '  -1.71%  ';
    "E46" = '  -0.52%  ';
    "E47" = '  -0.14%  ';
    "E48" = '  +3.18%  ';
    "E49" = '  -0.60%  ';
    "E50" = '  +2.30%  ';
    "E51" = '  +0.84%  ';
}

foreach ($addr in $otherValues.Keys) {
    $ws.Range($addr).Value = $otherValues[$addr]
}
